$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 138
$ws.Range("I6").Value = 140.4
$ws.Range("J6").Value = 134
$ws.Range("K6").Value = 421.2
$ws.Range("L6").Value = 402
$ws.Range("M6").Value = -309.2
$ws.Range("N6").Value = -626

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2625.6667
$ws.Range("I18").Value = 2825.5
$ws.Range("K18").Value = 2825.5
$ws.Range("M18").Value = -2541.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1417.2667
$ws.Range("I38").Value = 98.36364
$ws.Range("K38").Value = 295.09092
$ws.Range("M38").Value = 76.90908000000002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 83335864
$ws.Range("I40").Value = 2864
$ws.Range("J40").Value = 125002360
$ws.Range("K40").Value = 2864
$ws.Range("L40").Value = 125002360
$ws.Range("M40").Value = -2689
$ws.Range("N40").Value = -125002710

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4661.8184
$ws.Range("I100").Value = 2679.1667
$ws.Range("K100").Value = 2679.1667
$ws.Range("M100").Value = -2138.1667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 11894.292
$ws.Range("I116").Value = 7984.125
$ws.Range("K116").Value = 7984.125
$ws.Range("M116").Value = -4542.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1303891.4
$ws.Range("I137").Value = 1979.8
$ws.Range("J137").Value = 2605803
$ws.Range("K137").Value = 5939.4
$ws.Range("L137").Value = 7817409
$ws.Range("M137").Value = -3389.4
$ws.Range("N137").Value = -7822509

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3944.2812
$ws.Range("I138").Value = 1990.8518
$ws.Range("K138").Value = 5972.555399999999
$ws.Range("M138").Value = -832.5553999999993

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 95000
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 54999.5
$ws.Range("J44").Value = 89999
$ws.Range("L44").Value = 89999
$ws.Range("N44").Value = -90975

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 35002684
$ws.Range("I61").Value = 40003220
$ws.Range("K61").Value = 40003220
$ws.Range("M61").Value = -40003008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 927995.0600000001
$ws.Range("I74").Value = 1088254.6
$ws.Range("K74").Value = 1088254.6
$ws.Range("M74").Value = -1087380.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 927995.0600000001
$ws.Range("I77").Value = 1088254.6
$ws.Range("K77").Value = 5441273
$ws.Range("M77").Value = -5436905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5005723
$ws.Range("I132").Value = 6092.1665
$ws.Range("J132").Value = 50002400
$ws.Range("K132").Value = 18276.4995
$ws.Range("L132").Value = 150007200
$ws.Range("M132").Value = -15746.4995
$ws.Range("N132").Value = -150012260

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 35002684
$ws.Range("I136").Value = 40003220
$ws.Range("K136").Value = 120009660
$ws.Range("M136").Value = -120007110

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1279.7368
$ws.Range("J80").Value = 1183.4445
$ws.Range("L80").Value = 1183.4445
$ws.Range("N80").Value = -3179.4445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 1279.7368
$ws.Range("J83").Value = 1183.4445
$ws.Range("L83").Value = 5917.2225
$ws.Range("N83").Value = -15901.2225

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2954.818
$ws.Range("I99").Value = 2954.818
$ws.Range("K99").Value = 2954.818
$ws.Range("M99").Value = -1456.818

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 174999
$ws.Range("J132").Value = 174999
$ws.Range("L132").Value = 174999
$ws.Range("N132").Value = -185119

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5002721.5
$ws.Range("J134").Value = 8336475
$ws.Range("L134").Value = 25009425
$ws.Range("N134").Value = -25014495

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 50514140
$ws.Range("I31").Value = 62504388
$ws.Range("K31").Value = 62504388
$ws.Range("M31").Value = -62504093

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 50514140
$ws.Range("I34").Value = 62504388
$ws.Range("K34").Value = 62504388
$ws.Range("M34").Value = -62504186

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 99999
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 99999
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 99999
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -102371

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H90").Value = 99999
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 99999
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 299997
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -311853

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2127.825
$ws.Range("I107").Value = 1829.7742
$ws.Range("J107").Value = 3154.4443
$ws.Range("K107").Value = 1829.7742
$ws.Range("L107").Value = 3154.4443
$ws.Range("M107").Value = 90.22579999999994
$ws.Range("N107").Value = -6994.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2597.7778
$ws.Range("I132").Value = 2330.25
$ws.Range("J132").Value = 3132.8333
$ws.Range("K132").Value = 6990.75
$ws.Range("L132").Value = 9398.499899999999
$ws.Range("M132").Value = -4460.75
$ws.Range("N132").Value = -14458.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5863
$ws.Range("I134").Value = 5929.385
$ws.Range("K134").Value = 17788.155
$ws.Range("M134").Value = -15253.155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 129270.43
$ws.Range("J135").Value = 129270.43
$ws.Range("L135").Value = 129270.43
$ws.Range("N135").Value = -139410.43

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 13666.4
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 21444
$ws.Range("K22").Value = 6000
$ws.Range("L22").Value = 64332
$ws.Range("M22").Value = -5831
$ws.Range("N22").Value = -64670

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 13666.4
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 21444
$ws.Range("K27").Value = 6000
$ws.Range("L27").Value = 64332
$ws.Range("M27").Value = -5898
$ws.Range("N27").Value = -64536

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 11244
$ws.Range("I60").Value = 200
$ws.Range("M60").Value = -349

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5276.864
$ws.Range("I131").Value = 5006
$ws.Range("K131").Value = 15018
$ws.Range("M131").Value = -9978

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 635.61536
$ws.Range("J97").Value = 691.8333
$ws.Range("L97").Value = 691.8333
$ws.Range("N97").Value = -1683.8333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8303549.5
$ws.Range("I132").Value = 3347.5
$ws.Range("J132").Value = 27275440
$ws.Range("K132").Value = 10042.5
$ws.Range("L132").Value = 81826320
$ws.Range("M132").Value = -7512.5
$ws.Range("N132").Value = -81831380

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7064
$ws.Range("I7").Value = 6644.7144
$ws.Range("K7").Value = 6644.7144
$ws.Range("M7").Value = -6532.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3424.5789
$ws.Range("I40").Value = 3154.25
$ws.Range("K40").Value = 3154.25
$ws.Range("M40").Value = -3018.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5750
$ws.Range("J46").Value = 5750
$ws.Range("L46").Value = 5750
$ws.Range("N46").Value = -6126

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3410.4468
$ws.Range("I122").Value = 3315.7273
$ws.Range("J122").Value = 4799.6665
$ws.Range("K122").Value = 9947.1819
$ws.Range("L122").Value = 14398.9995
$ws.Range("M122").Value = -7497.1819
$ws.Range("N122").Value = -19298.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7064
$ws.Range("I126").Value = 6644.7144
$ws.Range("K126").Value = 19934.1432
$ws.Range("M126").Value = -17464.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5618.143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 149999.5
$ws.Range("I57").Value = 200000
$ws.Range("J57").Value = 99999
$ws.Range("K57").Value = 200000
$ws.Range("L57").Value = 99999
$ws.Range("M57").Value = -199246
$ws.Range("N57").Value = -101507

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1234.25
$ws.Range("J100").Value = 1149.75
$ws.Range("L100").Value = 2299.5
$ws.Range("N100").Value = -3381.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 839.5625
$ws.Range("I113").Value = 657.8
$ws.Range("J113").Value = 1142.5
$ws.Range("K113").Value = 1973.4
$ws.Range("L113").Value = 3427.5
$ws.Range("M113").Value = 196.6000000000001
$ws.Range("N113").Value = -7767.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3287.9333
$ws.Range("I122").Value = 2697.889
$ws.Range("K122").Value = 8093.667
$ws.Range("M122").Value = -5643.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 11473
$ws.Range("I126").Value = 11698.643
$ws.Range("J126").Value = 10841.2
$ws.Range("K126").Value = 35095.929
$ws.Range("L126").Value = 32523.6
$ws.Range("M126").Value = -32625.929
$ws.Range("N126").Value = -37463.60000000001
